$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Le Chat HS 300gr ADV23"
$ws.Range("A6").Value = 2958932

$ws.Range("B10").Value = "Le Chat Regular 2,5L ADV23"
$ws.Range("A10").Value = 2952074

$ws.Range("A10:B10").Select()
